# Swap the observation data between row 2 and row 3 on the active sheet.
# Only the columns whose values actually differ between the two rows are
# touched (A, B, E, F, G, H, Q, R, Z, AB); the remaining columns already
# hold identical values in both rows, so leaving them untouched keeps the
# rest of the sheet (including "empty" cells) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $cols) {
    $addrRow2 = $col + "2"
    $addrRow3 = $col + "3"

    $valueRow2 = $ws.Range($addrRow2).Value()
    $valueRow3 = $ws.Range($addrRow3).Value()

    $ws.Range($addrRow2).Value = $valueRow3
    $ws.Range($addrRow3).Value = $valueRow2
}
